$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-07 Saturday" "2026-02-08 Sunday"

Replace-Text "70×66=" "57×89="
Replace-Text "58×28=" "48×25="
Replace-Text "32×23=" "96×61="
Replace-Text "22×45=" "32×68="
Replace-Text "76×22=" "40×11="
Replace-Text "74×96=" "20×47="
Replace-Text "40×77=" "21×66="
Replace-Text "92×55=" "83×70="
Replace-Text "49×22=" "85×99="
Replace-Text "90×69=" "97×66="
Replace-Text "39×45=" "51×78="
Replace-Text "65×93=" "30×12="
Replace-Text "60×31=" "36×48="
Replace-Text "56×37=" "56×80="
Replace-Text "54×44=" "31×19="
Replace-Text "45×48=" "44×57="
Replace-Text "12×51=" "81×84="
Replace-Text "60×23=" "59×86="
Replace-Text "37×59=" "48×56="
Replace-Text "15×55=" "47×62="
Replace-Text "92×52=" "80×47="
Replace-Text "46×94=" "98×81="
Replace-Text "57×23=" "23×20="
Replace-Text "42×34=" "57×27="
Replace-Text "90×17=" "87×73="
